$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 40; this shifts the existing rows 40..115 down to 41..116
$ws.Rows(40).Insert()

# Populate the newly inserted row 40 with the new record
$ws.Range("A40").Value = 4
$ws.Range("B40").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C40").Value = "Los Lagos"
$ws.Range("D40").Value = 44469
$ws.Range("E40").Value = 10
$ws.Range("F40").Value = 100112039
$ws.Range("G40").Value = "Ciboulette"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 120
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = 3000
$ws.Range("N40").Value = "$/docena de atados"
$ws.Range("O40").Value = "Región Metropolitana"
$ws.Range("P40").Value = 1000
$ws.Range("Q40").Value = 3
$ws.Range("R40").Value = "Hortaliza"
